# Update the "想去人数" (interested-count) figures for the two sheets that
# carry the full event listing: "展览" (sheet 1) and "全部类型" (sheet 4).
# Sheets "演出" and "本地生活" only contain a header row, so nothing to do there.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value for column F (想去人数)
$updates = @{
    3 = 1248
    4 = 1500
    5 = 57
    6 = 6139
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
